$wb = $excel.ActiveWorkbook

$runManager = $wb.Worksheets.Item("RunManager")
$sheet1 = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# RunManager sheet: fix up Execute/Priority values for the existing test
# cases and append a brand-new "testcase3" row (row 6).
# ---------------------------------------------------------------------------
$runManager.Range("C3").Value = "yes"
$runManager.Range("E3").Value = "'1"

$runManager.Range("C4").Value = "yes"
$runManager.Range("E4").Value = "'1"

$runManager.Range("B5").Value = "To check whether the user can sort the apple laptops"
$runManager.Range("E5").Value = "'1"

$runManager.Range("A6").Value = "testcase3"
$runManager.Range("B6").Value = "Test Case 3"
$runManager.Range("C6").Value = "yes"
$runManager.Range("D6").Value = "'1"
$runManager.Range("E6").Value = "'1"

# ---------------------------------------------------------------------------
# Sheet1 ("RunManager" data-provider sheet): insert a new "browser" column
# (C) with the value "chrome" on every row, and append a new "testcase3"
# data row (row 11) for ELK / browser integration.
# ---------------------------------------------------------------------------
$sheet1.Columns.Item(3).Insert()
$sheet1.Columns.Item(3).ColumnWidth = $sheet1.Columns.Item(2).ColumnWidth

$sheet1.Range("C1").Value = "browser"
$sheet1.Range("C2").Value = "chrome"
$sheet1.Range("C3").Value = "chrome"
$sheet1.Range("C4").Value = "chrome"
$sheet1.Range("C5").Value = "chrome"
$sheet1.Range("C6").Value = "chrome"
$sheet1.Range("C7").Value = "chrome"
$sheet1.Range("C8").Value = "chrome"
$sheet1.Range("C9").Value = "chrome"
$sheet1.Range("C10").Value = "chrome"

# Row 9's trailing column used to be a blank (quote-prefixed) cell; give it
# real content while keeping the text formatting.
$sheet1.Range("F9").Value = "'jkn"

$sheet1.Range("A11").Value = "testcase3"
$sheet1.Range("B11").Value = "yes"
$sheet1.Range("C11").Value = "chrome"
$sheet1.Range("D11").Value = "sjdnc"
$sheet1.Range("E11").Value = "sdjn"
$sheet1.Range("F11").Value = "jkb"

# ---------------------------------------------------------------------------
# Selections / active tab: Sheet1 becomes the active tab with C3:C11
# selected; RunManager keeps a (no longer active) selection of E3:E5.
# ---------------------------------------------------------------------------
$runManager.Range("E3:E5").Select()

$sheet1.Activate()
$sheet1.Range("C3:C11").Select()
